# Raw and Clean Data from SSA for August 11th
# Applies the new day (2020-08-11, serial 44054) of data to the tracking workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper values
# ---------------------------------------------------------------------------
$newDate = 44054

# ---------------------------------------------------------------------------
# Sheet 1: out_vars - new row 73
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")

$ws1.Range("A73").Value = $newDate
$ws1.Range("A73").NumberFormat = "m/d/yy"
$ws1.Range("A73").WrapText = $true

$ws1.Range("B73").Value = 492522
$ws1.Range("C73").Value = 538333
$ws1.Range("D73").Value = 81259
$ws1.Range("E73").Value = 53929
$ws1.Range("F73").Value = 26.508460535773022
$ws1.Range("G73").Value = 130560
$ws1.Range("H73").Value = 10452
$ws1.Range("I73").Value = 12562
$ws1.Range("J73").Value = 1112114
$ws1.Range("B73:J73").WrapText = $true

# ---------------------------------------------------------------------------
# Sheet 2: dates_dx - new row 73
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")

$ws2.Range("A73").Value = $newDate
$ws2.Range("A73").NumberFormat = "m/d/yy"
$ws2.Range("A73").WrapText = $true

$ws2.Range("B73").Value = 0
$ws2.Range("C73").Value = 1
$ws2.Range("D73").Value = 0
$ws2.Range("E73").Value = 0
$ws2.Range("F73").Value = 1
$ws2.Range("G73").Value = 0
$ws2.Range("H73").Value = 0
$ws2.Range("I73").Value = 0
$ws2.Range("J73").Value = 0
$ws2.Range("K73").Value = 0
$ws2.Range("L73").Value = 4

for ($c = 2; $c -le 12; $c++) {
    $cell = $ws2.Cells.Item(73, $c)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# Sheet 3: dates_sx - new row 73
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")

$ws3.Range("A73").Value = $newDate
$ws3.Range("A73").NumberFormat = "m/d/yy"
$ws3.Range("A73").WrapText = $true

$ws3.Range("B73").Value = 0
$ws3.Range("C73").Value = 1
$ws3.Range("D73").Value = 0
$ws3.Range("E73").Value = 0
$ws3.Range("F73").Value = 0
$ws3.Range("G73").Value = 0
$ws3.Range("H73").Value = 1
$ws3.Range("I73").Value = 0
$ws3.Range("J73").Value = 0
$ws3.Range("K73").Value = 1
$ws3.Range("L73").Value = 0
$ws3.Range("M73").Value = 0
$ws3.Range("N73").Value = 0

for ($c = 2; $c -le 14; $c++) {
    $cell = $ws3.Cells.Item(73, $c)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# Sheet 4: dates_deaths - new row 73
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")

$ws4.Range("A73").Value = $newDate
$ws4.Range("A73").NumberFormat = "m/d/yy"
$ws4.Range("A73").WrapText = $true

$ws4.Range("B73").Value = 0
$ws4.Range("C73").Value = 0
$ws4.Range("D73").Value = 0
$ws4.Range("E73").Value = 0
$ws4.Range("F73").Value = 2
$ws4.Range("G73").Value = 1
$ws4.Range("H73").Value = 1
$ws4.Range("I73").Value = 1
$ws4.Range("J73").Value = 2

for ($c = 2; $c -le 10; $c++) {
    $cell = $ws4.Cells.Item(73, $c)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# Sheet 5: control_obs - fill the BU column (new date 2020-08-11) and totals
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")

$ws5.Range("BU1").Value = $newDate
$ws5.Range("BU2").Value = 5218
$ws5.Range("BU3").Value = 5017
$ws5.Range("BU4").Value = 5017
$ws5.Range("BU5").Value = 5017
$ws5.Range("BU6").Value = 5017
$ws5.Range("BU7").Value = 4327
$ws5.Range("BU8").Value = 6958

$ws5.Range("BU10").Value = 212
$ws5.Range("BU11").Value = 212
$ws5.Range("BU12").Value = 212
$ws5.Range("BU13").Value = 212
$ws5.Range("BU14").Value = 212
$ws5.Range("BU15").Value = 147
$ws5.Range("BU16").Value = 224

# Highlight the corrected observation (anomaly fix) in BT10:BU10
$ws5.Range("BT10:BU10").Interior.ThemeColor = 8

$ws5.Range("BU18").Value = 1213

# Extend the running-total row with the new column total
$ws5.Range("BU20").Formula = "=SUM(BU2:BU18)"
$ws5.Range("BU20").Borders.Item(7).LineStyle = 1
$ws5.Range("BU20").Borders.Item(10).LineStyle = 1
$ws5.Range("BU20").Borders.Item(8).LineStyle = 1
$ws5.Range("BU20").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# Sheet 7: anomalias - log the new entry for August 11th
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("anomalias")

$ws7.Range("A18").Value = "11 de agosto 2020"
$ws7.Range("B18").Value = "Las observaciones de sx_nal se mantuvieron"

# ---------------------------------------------------------------------------
# View state: make out_vars the active / selected sheet, restore selections
# ---------------------------------------------------------------------------
$ws2.Range("L73").Select()
$ws3.Range("O73").Select()
$ws4.Range("J73").Select()
$ws5.Range("BV16").Select()
$ws7.Range("B17").Select()

$ws1.Activate()
$ws1.Range("D64").Select()
